$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(121, 8).Value = 712
$ws.Cells.Item(121, 10).Value = 729.86664
$ws.Cells.Item(121, 12).Value = 2189.59992
$ws.Cells.Item(121, 14).Value = -5683.59992

$ws.Cells.Item(132, 8).Value = 1932098.9
$ws.Cells.Item(132, 9).Value = 1985698.8
$ws.Cells.Item(132, 10).Value = 2503
$ws.Cells.Item(132, 11).Value = 5957096.4
$ws.Cells.Item(132, 12).Value = 7509
$ws.Cells.Item(132, 13).Value = -5954566.4
$ws.Cells.Item(132, 14).Value = -12569

$ws.Cells.Item(138, 8).Value = 3224.5303
$ws.Cells.Item(138, 9).Value = 2472.7307
$ws.Cells.Item(138, 10).Value = 3713.2
$ws.Cells.Item(138, 11).Value = 7418.1921
$ws.Cells.Item(138, 12).Value = 11139.6
$ws.Cells.Item(138, 13).Value = -2278.1921
$ws.Cells.Item(138, 14).Value = -21419.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 23838.31
$ws.Cells.Item(32, 9).Value = 23277.59
$ws.Cells.Item(32, 11).Value = 23277.59
$ws.Cells.Item(32, 13).Value = -22990.59

$ws.Cells.Item(61, 8).Value = 1661.8572
$ws.Cells.Item(61, 9).Value = 1365.125
$ws.Cells.Item(61, 11).Value = 1365.125
$ws.Cells.Item(61, 13).Value = -1153.125

$ws.Cells.Item(122, 8).Value = 2030.4
$ws.Cells.Item(122, 9).Value = 2033.7778
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 6101.3334
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -3651.3334
$ws.Cells.Item(122, 14).Value = -10900

$ws.Cells.Item(132, 8).Value = 6763.5747
$ws.Cells.Item(132, 9).Value = 8237.031000000001
$ws.Cells.Item(132, 10).Value = 3620.2
$ws.Cells.Item(132, 11).Value = 24711.093
$ws.Cells.Item(132, 12).Value = 10860.6
$ws.Cells.Item(132, 13).Value = -22181.093
$ws.Cells.Item(132, 14).Value = -15920.6

$ws.Cells.Item(136, 8).Value = 1661.8572
$ws.Cells.Item(136, 9).Value = 1365.125
$ws.Cells.Item(136, 11).Value = 4095.375
$ws.Cells.Item(136, 13).Value = -1545.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 4222.107
$ws.Cells.Item(134, 9).Value = 5745.6206
$ws.Cells.Item(134, 10).Value = 2585.7407
$ws.Cells.Item(134, 11).Value = 17236.8618
$ws.Cells.Item(134, 12).Value = 7757.222099999999
$ws.Cells.Item(134, 13).Value = -14701.8618
$ws.Cells.Item(134, 14).Value = -12827.2221

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(41, 8).Value = 12050
$ws.Cells.Item(41, 10).Value = 19800
$ws.Cells.Item(41, 12).Value = 19800
$ws.Cells.Item(41, 14).Value = -20656

$ws.Cells.Item(50, 8).Value = 16650
$ws.Cells.Item(50, 10).Value = 16650
$ws.Cells.Item(50, 12).Value = 16650
$ws.Cells.Item(50, 14).Value = -17900

$ws.Cells.Item(51, 8).Value = 16133.333
$ws.Cells.Item(51, 10).Value = 22600
$ws.Cells.Item(51, 12).Value = 22600
$ws.Cells.Item(51, 14).Value = -24072

$ws.Cells.Item(58, 8).Value = 1426.2037
$ws.Cells.Item(58, 9).Value = 1287.575
$ws.Cells.Item(58, 10).Value = 1822.2858
$ws.Cells.Item(58, 11).Value = 1287.575
$ws.Cells.Item(58, 12).Value = 1822.2858
$ws.Cells.Item(58, 13).Value = -1084.575
$ws.Cells.Item(58, 14).Value = -2228.2858

$ws.Cells.Item(59, 8).Value = 36372.5
$ws.Cells.Item(59, 10).Value = 41663.332
$ws.Cells.Item(59, 12).Value = 41663.332
$ws.Cells.Item(59, 14).Value = -43953.332

$ws.Cells.Item(60, 8).Value = 21576.143
$ws.Cells.Item(60, 9).Value = 93
$ws.Cells.Item(60, 10).Value = 25156.666
$ws.Cells.Item(60, 11).Value = 93
$ws.Cells.Item(60, 12).Value = 25156.666
$ws.Cells.Item(60, 13).Value = 418
$ws.Cells.Item(60, 14).Value = -26178.666

$ws.Cells.Item(61, 8).Value = 16133.333
$ws.Cells.Item(61, 10).Value = 22600
$ws.Cells.Item(61, 12).Value = 22600
$ws.Cells.Item(61, 14).Value = -23296

$ws.Cells.Item(68, 8).Value = 42771.285
$ws.Cells.Item(68, 10).Value = 44899.832
$ws.Cells.Item(68, 12).Value = 44899.832
$ws.Cells.Item(68, 14).Value = -46397.832

$ws.Cells.Item(71, 8).Value = 42771.285
$ws.Cells.Item(71, 10).Value = 44899.832
$ws.Cells.Item(71, 12).Value = 134699.496
$ws.Cells.Item(71, 14).Value = -142187.496

$ws.Cells.Item(74, 8).Value = 24999.5
$ws.Cells.Item(74, 9).Value = 10000
$ws.Cells.Item(74, 10).Value = 39999
$ws.Cells.Item(74, 11).Value = 10000
$ws.Cells.Item(74, 12).Value = 39999
$ws.Cells.Item(74, 14).Value = -41747
$ws.Cells.Item(74, 13).Value = -9126

$ws.Cells.Item(77, 8).Value = 24999.5
$ws.Cells.Item(77, 9).Value = 10000
$ws.Cells.Item(77, 10).Value = 39999
$ws.Cells.Item(77, 11).Value = 30000
$ws.Cells.Item(77, 12).Value = 119997
$ws.Cells.Item(77, 14).Value = -128733
$ws.Cells.Item(77, 13).Value = -25632

$ws.Cells.Item(132, 8).Value = 3473808.2
$ws.Cells.Item(132, 9).Value = 1228.625
$ws.Cells.Item(132, 10).Value = 10418968
$ws.Cells.Item(132, 11).Value = 3685.875
$ws.Cells.Item(132, 12).Value = 31256904
$ws.Cells.Item(132, 13).Value = -1155.875
$ws.Cells.Item(132, 14).Value = -31261964

$ws.Cells.Item(134, 8).Value = 13726.75
$ws.Cells.Item(134, 9).Value = 13726.75
$ws.Cells.Item(134, 11).Value = 41180.25
$ws.Cells.Item(134, 13).Value = -38645.25

$ws.Cells.Item(136, 8).Value = 1426.2037
$ws.Cells.Item(136, 9).Value = 1287.575
$ws.Cells.Item(136, 10).Value = 1822.2858
$ws.Cells.Item(136, 11).Value = 3862.725
$ws.Cells.Item(136, 12).Value = 5466.857400000001
$ws.Cells.Item(136, 13).Value = -1312.725
$ws.Cells.Item(136, 14).Value = -10566.8574

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(107, 8).Value = 317.7
$ws.Cells.Item(107, 9).Value = 200.6
$ws.Cells.Item(107, 10).Value = 434.8
$ws.Cells.Item(107, 11).Value = 601.8
$ws.Cells.Item(107, 12).Value = 1304.4
$ws.Cells.Item(107, 13).Value = 1318.2
$ws.Cells.Item(107, 14).Value = -5144.4

$ws.Cells.Item(122, 8).Value = 1519.6
$ws.Cells.Item(122, 9).Value = 800
$ws.Cells.Item(122, 10).Value = 1999.3334
$ws.Cells.Item(122, 11).Value = 7200
$ws.Cells.Item(122, 12).Value = 17994.0006
$ws.Cells.Item(122, 13).Value = -4750
$ws.Cells.Item(122, 14).Value = -22894.0006

$ws.Cells.Item(129, 8).Value = 11906468
$ws.Cells.Item(129, 9).Value = 624.3333
$ws.Cells.Item(129, 10).Value = 20835850
$ws.Cells.Item(129, 11).Value = 1872.9999
$ws.Cells.Item(129, 12).Value = 62507550
$ws.Cells.Item(129, 13).Value = 3127.0001
$ws.Cells.Item(129, 14).Value = -62517550

$ws.Cells.Item(131, 8).Value = 2445.377
$ws.Cells.Item(131, 10).Value = 813.6949
$ws.Cells.Item(131, 12).Value = 2441.0847
$ws.Cells.Item(131, 14).Value = -12521.0847

$ws.Cells.Item(134, 8).Value = 1988.6364
$ws.Cells.Item(134, 9).Value = 1750
$ws.Cells.Item(134, 10).Value = 2500
$ws.Cells.Item(134, 11).Value = 5250
$ws.Cells.Item(134, 12).Value = 7500
$ws.Cells.Item(134, 13).Value = -180
$ws.Cells.Item(134, 14).Value = -17640

$ws.Cells.Item(136, 8).Value = 2010.3
$ws.Cells.Item(136, 9).Value = 1708.75
$ws.Cells.Item(136, 11).Value = 5126.25
$ws.Cells.Item(136, 13).Value = -26.25

$ws.Cells.Item(137, 8).Value = 46315000
$ws.Cells.Item(137, 9).Value = 27787702
$ws.Cells.Item(137, 10).Value = 58666530
$ws.Cells.Item(137, 11).Value = 83363106
$ws.Cells.Item(137, 12).Value = 175999590
$ws.Cells.Item(137, 13).Value = -83358006
$ws.Cells.Item(137, 14).Value = -176009790

$ws.Cells.Item(138, 8).Value = 2356.25
$ws.Cells.Item(138, 9).Value = 1443.1818
$ws.Cells.Item(138, 10).Value = 4365
$ws.Cells.Item(138, 11).Value = 4329.5454
$ws.Cells.Item(138, 12).Value = 13095
$ws.Cells.Item(138, 13).Value = 810.4546
$ws.Cells.Item(138, 14).Value = -23375

$ws.Cells.Item(139, 8).Value = 25001822
$ws.Cells.Item(139, 9).Value = 27779302
$ws.Cells.Item(139, 10).Value = 4500
$ws.Cells.Item(139, 11).Value = 83337906
$ws.Cells.Item(139, 12).Value = 13500
$ws.Cells.Item(139, 13).Value = -83332766
$ws.Cells.Item(139, 14).Value = -23780

$ws.Cells.Item(140, 8).Value = 2018.4615
$ws.Cells.Item(140, 9).Value = 901.7646999999999
$ws.Cells.Item(140, 11).Value = 2705.2941
$ws.Cells.Item(140, 13).Value = 2474.7059

$ws.Cells.Item(141, 8).Value = 4286.875
$ws.Cells.Item(141, 9).Value = 4432.222
$ws.Cells.Item(141, 10).Value = 4100
$ws.Cells.Item(141, 11).Value = 13296.666
$ws.Cells.Item(141, 12).Value = 12300
$ws.Cells.Item(141, 13).Value = -8116.665999999999
$ws.Cells.Item(141, 14).Value = -22660

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3002.4
$ws.Cells.Item(102, 9).Value = 2006
$ws.Cells.Item(102, 10).Value = 3666.6667
$ws.Cells.Item(102, 11).Value = 2006
$ws.Cells.Item(102, 12).Value = 3666.6667
$ws.Cells.Item(102, 13).Value = -384
$ws.Cells.Item(102, 14).Value = -6910.6667

$ws.Cells.Item(108, 8).Value = 57777
$ws.Cells.Item(108, 10).Value = 57777
$ws.Cells.Item(108, 12).Value = 57777
$ws.Cells.Item(108, 14).Value = -65457

$ws.Cells.Item(122, 8).Value = 20835482
$ws.Cells.Item(122, 9).Value = 31252010
$ws.Cells.Item(122, 10).Value = 2425
$ws.Cells.Item(122, 11).Value = 93756030
$ws.Cells.Item(122, 12).Value = 7275
$ws.Cells.Item(122, 13).Value = -93753580
$ws.Cells.Item(122, 14).Value = -12175

$ws.Cells.Item(132, 8).Value = 6159.3335
$ws.Cells.Item(132, 9).Value = 6930.476
$ws.Cells.Item(132, 10).Value = 3460.3333
$ws.Cells.Item(132, 11).Value = 20791.428
$ws.Cells.Item(132, 12).Value = 10380.9999
$ws.Cells.Item(132, 13).Value = -18261.428
$ws.Cells.Item(132, 14).Value = -15440.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 14785361
$ws.Cells.Item(107, 9).Value = 5435176
$ws.Cells.Item(107, 10).Value = 41667144
$ws.Cells.Item(107, 11).Value = 16305528
$ws.Cells.Item(107, 12).Value = 125001432
$ws.Cells.Item(107, 13).Value = -16303608
$ws.Cells.Item(107, 14).Value = -125005272

$ws.Cells.Item(113, 8).Value = 408.2
$ws.Cells.Item(113, 9).Value = 511.375
$ws.Cells.Item(113, 10).Value = 339.41666
$ws.Cells.Item(113, 11).Value = 1534.125
$ws.Cells.Item(113, 12).Value = 1018.24998
$ws.Cells.Item(113, 13).Value = 635.875
$ws.Cells.Item(113, 14).Value = -5358.24998

$ws.Cells.Item(114, 8).Value = 39750
$ws.Cells.Item(114, 10).Value = 39750
$ws.Cells.Item(114, 12).Value = 39750
$ws.Cells.Item(114, 14).Value = -48428
